# Profile sheet rework:
#  - set every "position" (column B) value to -1 (the engine now reports a
#    single sentinel instead of a 1-based ordinal)
#  - insert two new metric columns ("most_frequent_count", "csim") right
#    before "memory_consumed_bytes", shifting memory_consumed_bytes /
#    pattern_count / patterns one column to the right (T->V, U->W, V->X)
#  - populate the two new columns for every data row

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank columns at T:U - everything from T onward (memory_consumed_bytes,
# pattern_count, patterns) shifts right to V, W, X.
$ws.Range("T:U").Insert() | Out-Null

# New header labels for the freshly inserted columns.
$ws.Range("T1").Value = "most_frequent_count"
$ws.Range("U1").Value = "csim"

# Column B ("position") is now a constant sentinel for every attribute row.
$lastRow = 15
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 2).Value = -1
}

# Per-row values for the two new columns (most_frequent_count, csim).
$newValues = @{
    2  = @(8, 0)
    3  = @(21, 1)
    4  = @(12, 0)
    5  = @(17, 1)
    6  = @(2, 0)
    7  = @(5, 0)
    8  = @(2, 0)
    9  = @(4, 0)
    10 = @(24, 1)
    11 = @(24, 1)
    12 = @(10, 1)
    13 = @(4, 0)
    14 = @(4, 0.408)
    15 = @(4, 0)
}

foreach ($r in $newValues.Keys) {
    $pair = $newValues[$r]
    $ws.Range("T$r").Value = $pair[0]
    $ws.Range("U$r").Value = $pair[1]
}
